# Refresh currentAveragePrice / LevePrice / LeveProfit columns (H-N) for a batch of
# leve rows across all eight Masamune_Profits worksheets, matching the latest market
# snapshot pulled by the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 33
$ws.Range("H33").Value = 512.0909
$ws.Range("I33").Value = 114.666664
$ws.Range("K33").Value = 114.666664
$ws.Range("M33").Value = 114.333336

# row 39
$ws.Range("H39").Value = 111111330
$ws.Range("I39").Value = 121
$ws.Range("J39").Value = 333333730
$ws.Range("K39").Value = 363
$ws.Range("L39").Value = 1000001190
$ws.Range("M39").Value = -67
$ws.Range("N39").Value = -1000001782

# row 44
$ws.Range("H44").Value = 35272.375
$ws.Range("J44").Value = 35272.375
$ws.Range("L44").Value = 35272.375
$ws.Range("N44").Value = -36196.375

# row 64
$ws.Range("H64").Value = 30305460
$ws.Range("I64").Value = 166666670
$ws.Range("J64").Value = 7578590.5
$ws.Range("K64").Value = 166666670
$ws.Range("L64").Value = 7578590.5
$ws.Range("M64").Value = -166666422
$ws.Range("N64").Value = -7579086.5

# row 67
$ws.Range("H67").Value = 30305460
$ws.Range("I67").Value = 166666670
$ws.Range("J67").Value = 7578590.5
$ws.Range("K67").Value = 166666670
$ws.Range("L67").Value = 7578590.5
$ws.Range("M67").Value = -166665812
$ws.Range("N67").Value = -7580306.5

# row 70
$ws.Range("H70").Value = 1852.9474
$ws.Range("I70").Value = 2000
$ws.Range("J70").Value = 1813.7333
$ws.Range("K70").Value = 6000
$ws.Range("L70").Value = 5441.199900000001
$ws.Range("M70").Value = -5730
$ws.Range("N70").Value = -5981.199900000001

# row 73
$ws.Range("H73").Value = 1852.9474
$ws.Range("I73").Value = 2000
$ws.Range("J73").Value = 1813.7333
$ws.Range("K73").Value = 6000
$ws.Range("L73").Value = 5441.199900000001
$ws.Range("M73").Value = -5064
$ws.Range("N73").Value = -7313.199900000001

# row 86
$ws.Range("H86").Value = 3847341.2
$ws.Range("I86").Value = 5001160.5
$ws.Range("J86").Value = 1277.1666
$ws.Range("K86").Value = 5001160.5
$ws.Range("L86").Value = 1277.1666
$ws.Range("M86").Value = -5000037.5
$ws.Range("N86").Value = -3523.1666

# row 89
$ws.Range("H89").Value = 3847341.2
$ws.Range("I89").Value = 5001160.5
$ws.Range("J89").Value = 1277.1666
$ws.Range("K89").Value = 25005802.5
$ws.Range("L89").Value = 6385.833000000001
$ws.Range("M89").Value = -25000186.5
$ws.Range("N89").Value = -17617.833

# row 127
$ws.Range("H127").Value = 2051.8333
$ws.Range("I127").Value = 1024.75
$ws.Range("J127").Value = 2345.2856
$ws.Range("K127").Value = 3074.25
$ws.Range("L127").Value = 7035.8568
$ws.Range("M127").Value = 1885.75
$ws.Range("N127").Value = -16955.8568

$ws = $wb.Worksheets.Item("ARM")
# row 2
$ws.Range("H2").Value = 3280
$ws.Range("J2").Value = 10000
$ws.Range("L2").Value = 10000
$ws.Range("N2").Value = -10226

# row 4
$ws.Range("H4").Value = 466.66666
$ws.Range("I4").Value = 300
$ws.Range("J4").Value = 550
$ws.Range("K4").Value = 300
$ws.Range("L4").Value = 550
$ws.Range("M4").Value = -184
$ws.Range("N4").Value = -782

# row 116
$ws.Range("H116").Value = 3280
$ws.Range("J116").Value = 10000
$ws.Range("L116").Value = 10000
$ws.Range("N116").Value = -14588

# row 132
$ws.Range("H132").Value = 8774007
$ws.Range("I132").Value = 12501315
$ws.Range("J132").Value = 3871.5881
$ws.Range("K132").Value = 37503945
$ws.Range("L132").Value = 11614.7643
$ws.Range("M132").Value = -37501415
$ws.Range("N132").Value = -16674.7643

$ws = $wb.Worksheets.Item("BSM")
# row 3
$ws.Range("H3").Value = 3280
$ws.Range("J3").Value = 10000
$ws.Range("L3").Value = 10000
$ws.Range("N3").Value = -10228

# row 134
$ws.Range("H134").Value = 1592.4546
$ws.Range("I134").Value = 1227.4073
$ws.Range("J134").Value = 3235.1667
$ws.Range("K134").Value = 3682.2219
$ws.Range("L134").Value = 9705.500100000001
$ws.Range("M134").Value = -1147.2219
$ws.Range("N134").Value = -14775.5001

$ws = $wb.Worksheets.Item("CRP")
# row 62
$ws.Range("H62").Value = 3612.4167
$ws.Range("I62").Value = 4190
$ws.Range("J62").Value = 3199.8572
$ws.Range("K62").Value = 4190
$ws.Range("L62").Value = 3199.8572
$ws.Range("M62").Value = -3566
$ws.Range("N62").Value = -4447.8572

# row 65
$ws.Range("H65").Value = 3612.4167
$ws.Range("I65").Value = 4190
$ws.Range("J65").Value = 3199.8572
$ws.Range("K65").Value = 20950
$ws.Range("L65").Value = 15999.286
$ws.Range("M65").Value = -17830
$ws.Range("N65").Value = -22239.286

# row 132
$ws.Range("H132").Value = 32435.588
$ws.Range("I132").Value = 1618.909
$ws.Range("J132").Value = 110662.54
$ws.Range("K132").Value = 4856.727000000001
$ws.Range("L132").Value = 331987.62
$ws.Range("M132").Value = -2326.727000000001
$ws.Range("N132").Value = -337047.62

$ws = $wb.Worksheets.Item("CUL")
# row 117
$ws.Range("H117").Value = 1638.9445
$ws.Range("I117").Value = 427.83334
$ws.Range("J117").Value = 2244.5
$ws.Range("K117").Value = 1283.50002
$ws.Range("L117").Value = 6733.5
$ws.Range("M117").Value = 2158.49998
$ws.Range("N117").Value = -13617.5

# row 139
$ws.Range("H139").Value = 117084.16
$ws.Range("I139").Value = 117084.16
$ws.Range("K139").Value = 351252.48
$ws.Range("M139").Value = -346112.48

$ws = $wb.Worksheets.Item("GSM")
# row 64
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

# row 67
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

# row 132
$ws.Range("H132").Value = 2385.3489
$ws.Range("I132").Value = 1639.4
$ws.Range("J132").Value = 4106.769
$ws.Range("K132").Value = 4918.200000000001
$ws.Range("L132").Value = 12320.307
$ws.Range("M132").Value = -2388.200000000001
$ws.Range("N132").Value = -17380.307

$ws = $wb.Worksheets.Item("LTW")
# row 20
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()

# row 21
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()

# row 34
$ws.Range("H34").Value = 21964.5
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 21964.5
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 21964.5
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -22308.5

# row 38
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("N38").ClearContents()

# row 39
$ws.Range("H39").Value = 34215
$ws.Range("I39").Value = 20000
$ws.Range("J39").Value = 38953.332
$ws.Range("K39").Value = 20000
$ws.Range("L39").Value = 38953.332
$ws.Range("M39").Value = -19540
$ws.Range("N39").Value = -39873.332

# row 40
$ws.Range("H40").Value = 5584.3335
$ws.Range("I40").Value = 5375.625
$ws.Range("J40").Value = 6001.75
$ws.Range("K40").Value = 5375.625
$ws.Range("L40").Value = 6001.75
$ws.Range("M40").Value = -5239.625
$ws.Range("N40").Value = -6273.75

# row 45
$ws.Range("H45").Value = 18217.5
$ws.Range("I45").Value = 19000
$ws.Range("J45").Value = 18105.715
$ws.Range("K45").Value = 19000
$ws.Range("L45").Value = 18105.715
$ws.Range("M45").Value = -18593
$ws.Range("N45").Value = -18919.715

# row 46
$ws.Range("H46").Value = 4150.115
$ws.Range("I46").Value = 853.4
$ws.Range("J46").Value = 8645.637000000001
$ws.Range("K46").Value = 853.4
$ws.Range("L46").Value = 8645.637000000001
$ws.Range("M46").Value = -665.4
$ws.Range("N46").Value = -9021.637000000001

# row 51
$ws.Range("H51").Value = 32000
$ws.Range("J51").Value = 32000
$ws.Range("L51").Value = 32000
$ws.Range("N51").Value = -32956

# row 56
$ws.Range("H56").Value = 23053.5
$ws.Range("J56").Value = 23053.5
$ws.Range("L56").Value = 23053.5
$ws.Range("N56").Value = -24435.5

# row 57
$ws.Range("H57").Value = 68533
$ws.Range("J57").Value = 68533
$ws.Range("L57").Value = 68533
$ws.Range("N57").Value = -69665

# row 58
$ws.Range("H58").Value = 10103
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 10103
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 10103
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -10623

# row 68
$ws.Range("H68").Value = 3776.12
$ws.Range("I68").Value = 3729.5881
$ws.Range("J68").Value = 3875
$ws.Range("K68").Value = 3729.5881
$ws.Range("L68").Value = 3875
$ws.Range("M68").Value = -2980.5881
$ws.Range("N68").Value = -5373

# row 71
$ws.Range("H71").Value = 3776.12
$ws.Range("I71").Value = 3729.5881
$ws.Range("J71").Value = 3875
$ws.Range("K71").Value = 18647.9405
$ws.Range("L71").Value = 19375
$ws.Range("M71").Value = -14903.9405
$ws.Range("N71").Value = -26863

# row 100
$ws.Range("H100").Value = 3456.1538
$ws.Range("I100").Value = 2488.3333
$ws.Range("J100").Value = 4285.7144
$ws.Range("K100").Value = 2488.3333
$ws.Range("L100").Value = 4285.7144
$ws.Range("M100").Value = -1947.3333
$ws.Range("N100").Value = -5367.7144

# row 122
$ws.Range("H122").Value = 61062
$ws.Range("I122").Value = 73668.14
$ws.Range("J122").Value = 2233.3333
$ws.Range("K122").Value = 221004.42
$ws.Range("L122").Value = 6699.999899999999
$ws.Range("M122").Value = -218554.42
$ws.Range("N122").Value = -11599.9999

# row 132
$ws.Range("H132").Value = 4048.3333
$ws.Range("J132").Value = 4828.6924
$ws.Range("L132").Value = 14486.0772
$ws.Range("N132").Value = -19546.0772

$ws = $wb.Worksheets.Item("WVR")
# row 62
$ws.Range("H62").Value = 2999.9473
$ws.Range("I62").Value = 2999
$ws.Range("K62").Value = 2999
$ws.Range("M62").Value = -2375

# row 65
$ws.Range("H65").Value = 2999.9473
$ws.Range("I65").Value = 2999
$ws.Range("K65").Value = 14995
$ws.Range("M65").Value = -11875

# row 96
$ws.Range("H96").Value = 2340.6
$ws.Range("I96").Value = 2450.75
$ws.Range("J96").Value = 1900
$ws.Range("K96").Value = 2450.75
$ws.Range("L96").Value = 1900
$ws.Range("M96").Value = -1077.75
$ws.Range("N96").Value = -4646

# row 122
$ws.Range("H122").Value = 815.5
$ws.Range("I122").Value = 815.5
$ws.Range("K122").Value = 2446.5
$ws.Range("M122").Value = 3.5

# row 132
$ws.Range("H132").Value = 1714.122
$ws.Range("J132").Value = 2062.2273
$ws.Range("L132").Value = 6186.6819
$ws.Range("N132").Value = -11246.6819
